# Update crypto price/volume data per GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.036.07"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.650.82"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5285"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06303"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.483"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.657.76"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5465"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "0.0₅8102"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "26.050.40"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.568"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "139.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1244"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.278"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.412"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05948"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.276"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.501"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.249"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.540"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.413"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.757"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5663"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.866"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8475"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "1.007.58"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("D44").Value = "1.789.10"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.485"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4286"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05150"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.835"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.04%  "
